$d = $word.ActiveDocument

# Replace every occurrence of $oldText with $newText using ordinary Find/Replace
# (wdReplaceOne in a loop so each hit is handled individually). Word's Find engine
# "smartens" plain straight apostrophes into curly ones as part of the replace,
# even though the source text we supply only contains straight apostrophes; after
# each hit we compare the resulting range text against the intended literal text
# and, if AutoCorrect altered it, overwrite just that (already-isolated) range via
# a direct Range.Text assignment, which does not invoke the smart-quote engine.
function Replace-AllText($doc, $oldText, $newText) {
    $count = 0
    while ($true) {
        $rng = $doc.Content
        $rng.Start = 0
        $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)
        if (-not $found) { break }
        if ($rng.Text -cne $newText) {
            $rng.Text = $newText
        }
        $count = $count + 1
        if ($count -gt 100) { break }
    }
    return $count
}

# 1. Title heading + bold closing line (identical text, two occurrences)
Replace-AllText $d "Play Alchemist's Gold Slot for Free - Review & Features" "Play Alchemist's Gold for Free"

# 2. "What we like" bullet points
Replace-AllText $d "Mystery Symbol function for increased chances of winning" "Simple gameplay with a Mystery Symbol function"
Replace-AllText $d "Attractive and well-crafted graphics with theme-related symbols" "Well-crafted graphics that stand out"
Replace-AllText $d "Standard structure with 5 reels and 10 lines" "Symbols that reflect the alchemy theme"

# 3. "What we don't like" bullet points
Replace-AllText $d "Bonus feature is inconsistent" "Inconsistent bonus feature"
Replace-AllText $d "50/50 game is risky for players" "Risky 50/50 game"

# 4. Italic summary line near the end
Replace-AllText $d "Discover the alchemy-themed slot game Alchemist's Gold. Play for free and experience the Mystery Symbol function and chance to double winnings." "Read our review of Alchemist's Gold and play for free. Discover its features and gameplay."
